$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: reuse existing "2025-04-10" value for period_start/period_end,
# add a new event description for the entrance ceremony day.
$ws.Range("A6").Value = "2025-04-10"
$ws.Range("B6").Value = "2025-04-10"
$ws.Range("F6").Value = "4月10日の予定"

# Row 7: new date (2025-04-01) with its own event description.
$ws.Range("A7").Value = "2025-04-01"
$ws.Range("B7").Value = "2025-04-01"
$ws.Range("F7").Value = "入学式"

# Update the active selection to match the edited range.
[void]$ws.Range("B6:B7").Select()
